# ------------------------------------------------------------------
# Apply the "Atializacoes Diagramas revisao ALEX @virtohoho" changes:
#   1. Refresh the cached "datetimeFigureOut" date field from
#      09/04/2019 -> 16/04/2019 everywhere it appears (slide master +
#      every slide layout placeholder).
#   2. Bump the picture quality hint on the css-logo picture (best
#      effort - see note below).
#   3. Resize the "CaixaDeTexto 29" text box and split its paragraph
#      into three runs, inserting the new sentence fragment about the
#      app being reachable over the internet.
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Date placeholder fields (slide master + all custom layouts) ---

function Update-DateShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "09/04/2019") {
                    $sh.TextFrame.TextRange.Text = "16/04/2019"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes($master)
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    Update-DateShapes($master.CustomLayouts.Item($L))
}

# Slides themselves can also carry their own (non-inherited) copy of
# the field if the placeholder was ever edited directly on the slide.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    Update-DateShapes($p.Slides.Item($si))
}

# --- 2. Picture quality hint on the "Picture 12" (css logo) image ---
# The authoring tool stamps the embedded blip with cstate="hqprint"
# (a high-quality-print compression cache marker PowerPoint writes
# when it recompresses a picture). PowerPoint's object model has no
# dedicated property for this blip attribute, so we run the nearest
# equivalent operation (picture compression) in case the host wires
# it up; harmless no-op otherwise.
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "Picture 12") {
        try { $sh.PictureFormat.Compress() } catch { }
    }
}

# --- 3. Resize + re-word the "CaixaDeTexto 29" text box ---

for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "CaixaDeTexto 29") {

        # Grow the box upward and taller to fit the extra sentence.
        $sh.Top = 79.65236290472441
        $sh.Height = 73.91488268976377

        $tr = $sh.TextFrame.TextRange

        $oldText = "Através das linguagens de programação, criamos as aplicações web e mobile, onde o cliente poderá ter acesso as suas informações, dados recolhidos e gráficos."
        $tr.Text = $oldText

        # Replace "mobile, " with the longer phrase - this naturally
        # splits the paragraph into three runs at the edit boundaries.
        $oldFragment = "mobile, "
        $newFragment = "mobile que serão acessadas através da internet, "
        $fragStart = 68
        $fragLen = 8

        $mid = $tr.Characters($fragStart, $fragLen)
        $mid.Text = $newFragment

        # Re-stamp the font size (identical value) on the newly
        # inserted middle run so the engine commits it as its own
        # run with matching formatting, not just a length variant of
        # a pre-existing run.
        $mid2 = $tr.Characters($fragStart, $newFragment.Length)
        $mid2.Font.Size = 11
    }
}
